# Auto-generated Excel COM-interop script
# Applies numeric corrections to columns H-N across multiple sheets
# as described in the commit's canonical OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 195.95454
$ws.Cells.Item(33, 9).Value = 92.388885
$ws.Cells.Item(33, 11).Value = 92.388885
$ws.Cells.Item(33, 13).Value = 136.611115
$ws.Cells.Item(53, 8).Value = 863.8276
$ws.Cells.Item(53, 9).Value = 1124.2273
$ws.Cells.Item(53, 10).Value = 45.42857
$ws.Cells.Item(53, 11).Value = 1124.2273
$ws.Cells.Item(53, 12).Value = 45.42857
$ws.Cells.Item(53, 13).Value = -487.2273
$ws.Cells.Item(53, 14).Value = -1319.42857
$ws.Cells.Item(76, 8).Value = 13333.333
$ws.Cells.Item(76, 9).Value = 10000
$ws.Cells.Item(76, 10).Value = 20000
$ws.Cells.Item(76, 11).Value = 10000
$ws.Cells.Item(76, 12).Value = 20000
$ws.Cells.Item(76, 13).Value = -9685
$ws.Cells.Item(76, 14).Value = -20630
$ws.Cells.Item(79, 8).Value = 13333.333
$ws.Cells.Item(79, 9).Value = 10000
$ws.Cells.Item(79, 10).Value = 20000
$ws.Cells.Item(79, 11).Value = 10000
$ws.Cells.Item(79, 12).Value = 20000
$ws.Cells.Item(79, 13).Value = -8908
$ws.Cells.Item(79, 14).Value = -22184
$ws.Cells.Item(86, 8).Value = 4387
$ws.Cells.Item(86, 9).Value = 4125.7144
$ws.Cells.Item(86, 10).Value = 4569.9
$ws.Cells.Item(86, 11).Value = 4125.7144
$ws.Cells.Item(86, 12).Value = 4569.9
$ws.Cells.Item(86, 13).Value = -3002.7144
$ws.Cells.Item(86, 14).Value = -6815.9
$ws.Cells.Item(89, 8).Value = 4387
$ws.Cells.Item(89, 9).Value = 4125.7144
$ws.Cells.Item(89, 10).Value = 4569.9
$ws.Cells.Item(89, 11).Value = 20628.572
$ws.Cells.Item(89, 12).Value = 22849.5
$ws.Cells.Item(89, 13).Value = -15012.572
$ws.Cells.Item(89, 14).Value = -34081.5
$ws.Cells.Item(98, 8).Value = 6189.864
$ws.Cells.Item(98, 9).Value = 6189.864
$ws.Cells.Item(98, 11).Value = 6189.864
$ws.Cells.Item(98, 13).Value = -4691.864
$ws.Cells.Item(106, 8).Value = 3856.2727
$ws.Cells.Item(106, 9).Value = 3641.3
$ws.Cells.Item(106, 11).Value = 3641.3
$ws.Cells.Item(106, 13).Value = -3010.3
$ws.Cells.Item(112, 8).Value = 2612.16
$ws.Cells.Item(112, 9).Value = 750
$ws.Cells.Item(112, 10).Value = 3200.2104
$ws.Cells.Item(112, 11).Value = 2250
$ws.Cells.Item(112, 12).Value = 9600.6312
$ws.Cells.Item(112, 13).Value = -1142
$ws.Cells.Item(112, 14).Value = -11816.6312
$ws.Cells.Item(122, 8).Value = 6189.864
$ws.Cells.Item(122, 9).Value = 6189.864
$ws.Cells.Item(122, 11).Value = 18569.592
$ws.Cells.Item(122, 13).Value = -16119.592
$ws.Cells.Item(132, 8).Value = 5468237.5
$ws.Cells.Item(132, 9).Value = 7577185.5
$ws.Cells.Item(132, 11).Value = 22731556.5
$ws.Cells.Item(132, 13).Value = -22729026.5
$ws.Cells.Item(135, 8).Value = 956.6667
$ws.Cells.Item(135, 9).Value = 377.42856
$ws.Cells.Item(135, 11).Value = 3396.85704
$ws.Cells.Item(135, 13).Value = -861.8570399999999
$ws.Cells.Item(137, 8).Value = 1069.4667
$ws.Cells.Item(137, 9).Value = 860.8570999999999
$ws.Cells.Item(137, 10).Value = 1334.9697
$ws.Cells.Item(137, 11).Value = 2582.5713
$ws.Cells.Item(137, 12).Value = 4004.9091
$ws.Cells.Item(137, 13).Value = -32.57129999999961
$ws.Cells.Item(137, 14).Value = -9104.909100000001
$ws.Cells.Item(138, 8).Value = 1352.36
$ws.Cells.Item(138, 9).Value = 608.2222
$ws.Cells.Item(138, 10).Value = 1627.589
$ws.Cells.Item(138, 11).Value = 1824.6666
$ws.Cells.Item(138, 12).Value = 4882.767
$ws.Cells.Item(138, 13).Value = 3315.3334
$ws.Cells.Item(138, 14).Value = -15162.767
$ws.Cells.Item(141, 8).Value = 656.4138
$ws.Cells.Item(141, 9).Value = 572.7143
$ws.Cells.Item(141, 10).Value = 3000
$ws.Cells.Item(141, 11).Value = 1718.1429
$ws.Cells.Item(141, 12).Value = 9000
$ws.Cells.Item(141, 13).Value = 3461.8571
$ws.Cells.Item(141, 14).Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6010.4375
$ws.Cells.Item(32, 9).Value = 5661.1333
$ws.Cells.Item(32, 10).Value = 11250
$ws.Cells.Item(32, 11).Value = 5661.1333
$ws.Cells.Item(32, 12).Value = 11250
$ws.Cells.Item(32, 13).Value = -5374.1333
$ws.Cells.Item(32, 14).Value = -11824
$ws.Cells.Item(45, 8).Value = 1373.6842
$ws.Cells.Item(45, 9).Value = 1386.7059
$ws.Cells.Item(45, 10).Value = 1263
$ws.Cells.Item(45, 11).Value = 1386.7059
$ws.Cells.Item(45, 12).Value = 1263
$ws.Cells.Item(45, 13).Value = -1009.7059
$ws.Cells.Item(45, 14).Value = -2017
$ws.Cells.Item(88, 8).Value = 2497
$ws.Cells.Item(88, 9).Value = 2100
$ws.Cells.Item(88, 10).Value = 2536.7
$ws.Cells.Item(88, 11).Value = 2100
$ws.Cells.Item(88, 12).Value = 2536.7
$ws.Cells.Item(88, 13).Value = -1694
$ws.Cells.Item(88, 14).Value = -3348.7
$ws.Cells.Item(91, 8).Value = 2497
$ws.Cells.Item(91, 9).Value = 2100
$ws.Cells.Item(91, 10).Value = 2536.7
$ws.Cells.Item(91, 11).Value = 2100
$ws.Cells.Item(91, 12).Value = 2536.7
$ws.Cells.Item(91, 13).Value = -696
$ws.Cells.Item(91, 14).Value = -5344.7
$ws.Cells.Item(97, 8).Value = 395.91306
$ws.Cells.Item(97, 9).Value = 335.3
$ws.Cells.Item(97, 10).Value = 800
$ws.Cells.Item(97, 11).Value = 335.3
$ws.Cells.Item(97, 12).Value = 800
$ws.Cells.Item(97, 13).Value = 160.7
$ws.Cells.Item(97, 14).Value = -1792
$ws.Cells.Item(110, 8).Value = 1557.8636
$ws.Cells.Item(110, 9).Value = 1203.5294
$ws.Cells.Item(110, 11).Value = 1203.5294
$ws.Cells.Item(110, 13).Value = 841.4706000000001
$ws.Cells.Item(122, 8).Value = 2664.8572
$ws.Cells.Item(122, 9).Value = 2553.5
$ws.Cells.Item(122, 10).Value = 3333
$ws.Cells.Item(122, 11).Value = 7660.5
$ws.Cells.Item(122, 12).Value = 9999
$ws.Cells.Item(122, 13).Value = -5210.5
$ws.Cells.Item(122, 14).Value = -14899
$ws.Cells.Item(132, 8).Value = 3510.8635
$ws.Cells.Item(132, 9).Value = 5103.778
$ws.Cells.Item(132, 11).Value = 15311.334
$ws.Cells.Item(132, 13).Value = -12781.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 14707172
$ws.Cells.Item(94, 9).Value = 19231986
$ws.Cells.Item(94, 10).Value = 1523.25
$ws.Cells.Item(94, 11).Value = 19231986
$ws.Cells.Item(94, 12).Value = 1523.25
$ws.Cells.Item(94, 13).Value = -19231535
$ws.Cells.Item(94, 14).Value = -2425.25
$ws.Cells.Item(134, 8).Value = 4254.511
$ws.Cells.Item(134, 9).Value = 1464.2727
$ws.Cells.Item(134, 10).Value = 11927.667
$ws.Cells.Item(134, 11).Value = 4392.8181
$ws.Cells.Item(134, 12).Value = 35783.001
$ws.Cells.Item(134, 13).Value = -1857.8181
$ws.Cells.Item(134, 14).Value = -40853.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2278.1
$ws.Cells.Item(31, 9).Value = 2240.1052
$ws.Cells.Item(31, 10).Value = 3000
$ws.Cells.Item(31, 11).Value = 2240.1052
$ws.Cells.Item(31, 12).Value = 3000
$ws.Cells.Item(31, 13).Value = -1945.1052
$ws.Cells.Item(31, 14).Value = -3590
$ws.Cells.Item(34, 8).Value = 2278.1
$ws.Cells.Item(34, 9).Value = 2240.1052
$ws.Cells.Item(34, 10).Value = 3000
$ws.Cells.Item(34, 11).Value = 2240.1052
$ws.Cells.Item(34, 12).Value = 3000
$ws.Cells.Item(34, 13).Value = -2038.1052
$ws.Cells.Item(34, 14).Value = -3404
$ws.Cells.Item(132, 8).Value = 5141.7
$ws.Cells.Item(132, 9).Value = 5202.0386
$ws.Cells.Item(132, 11).Value = 15606.1158
$ws.Cells.Item(132, 13).Value = -13076.1158
$ws.Cells.Item(134, 8).Value = 831.6984
$ws.Cells.Item(134, 9).Value = 807.2632
$ws.Cells.Item(134, 11).Value = 2421.7896
$ws.Cells.Item(134, 13).Value = 113.2103999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 461
$ws.Cells.Item(50, 9).Value = 105
$ws.Cells.Item(50, 10).Value = 550
$ws.Cells.Item(50, 11).Value = 315
$ws.Cells.Item(50, 12).Value = 1650
$ws.Cells.Item(50, 13).Value = 166
$ws.Cells.Item(50, 14).Value = -2612
$ws.Cells.Item(53, 8).Value = 461
$ws.Cells.Item(53, 9).Value = 105
$ws.Cells.Item(53, 10).Value = 550
$ws.Cells.Item(53, 11).Value = 315
$ws.Cells.Item(53, 12).Value = 1650
$ws.Cells.Item(53, 13).Value = 166
$ws.Cells.Item(53, 14).Value = -2612
$ws.Cells.Item(131, 8).Value = 33334850
$ws.Cells.Item(131, 10).Value = 2027.5238
$ws.Cells.Item(131, 12).Value = 6082.5714
$ws.Cells.Item(131, 14).Value = -16162.5714
$ws.Cells.Item(137, 8).Value = 24195842
$ws.Cells.Item(137, 9).Value = 57693216
$ws.Cells.Item(137, 10).Value = 3295
$ws.Cells.Item(137, 11).Value = 173079648
$ws.Cells.Item(137, 12).Value = 9885
$ws.Cells.Item(137, 13).Value = -173074548
$ws.Cells.Item(137, 14).Value = -20085

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(45, 8).Value = 34428.4
$ws.Cells.Item(45, 10).Value = 34428.4
$ws.Cells.Item(45, 12).Value = 34428.4
$ws.Cells.Item(45, 14).Value = -35546.4
$ws.Cells.Item(80, 8).Value = 2273.3333
$ws.Cells.Item(80, 9).Value = 1518.1818
$ws.Cells.Item(80, 10).Value = 4350
$ws.Cells.Item(80, 11).Value = 1518.1818
$ws.Cells.Item(80, 12).Value = 4350
$ws.Cells.Item(80, 13).Value = -520.1818000000001
$ws.Cells.Item(80, 14).Value = -6346
$ws.Cells.Item(83, 8).Value = 2273.3333
$ws.Cells.Item(83, 9).Value = 1518.1818
$ws.Cells.Item(83, 10).Value = 4350
$ws.Cells.Item(83, 11).Value = 7590.909000000001
$ws.Cells.Item(83, 12).Value = 21750
$ws.Cells.Item(83, 13).Value = -2598.909000000001
$ws.Cells.Item(83, 14).Value = -31734
$ws.Cells.Item(97, 8).Value = 607.1429000000001
$ws.Cells.Item(97, 9).Value = 629.9231
$ws.Cells.Item(97, 11).Value = 629.9231
$ws.Cells.Item(97, 13).Value = -133.9231
$ws.Cells.Item(113, 8).Value = 1797.3572
$ws.Cells.Item(113, 9).Value = 1776.0834
$ws.Cells.Item(113, 10).Value = 1925
$ws.Cells.Item(113, 11).Value = 1776.0834
$ws.Cells.Item(113, 12).Value = 1925
$ws.Cells.Item(113, 13).Value = 393.9166
$ws.Cells.Item(113, 14).Value = -6265
$ws.Cells.Item(126, 8).Value = 2139.25
$ws.Cells.Item(126, 9).Value = 1860
$ws.Cells.Item(126, 10).Value = 2604.6667
$ws.Cells.Item(126, 11).Value = 5580
$ws.Cells.Item(126, 12).Value = 7814.000100000001
$ws.Cells.Item(126, 13).Value = -3110
$ws.Cells.Item(126, 14).Value = -12754.0001
$ws.Cells.Item(135, 8).Value = 42250
$ws.Cells.Item(135, 10).Value = 34500
$ws.Cells.Item(135, 12).Value = 34500
$ws.Cells.Item(135, 14).Value = -44640

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94, 8).Value = 19995.9
$ws.Cells.Item(94, 10).Value = 19995.9
$ws.Cells.Item(94, 12).Value = 19995.9
$ws.Cells.Item(94, 14).Value = -21347.9
$ws.Cells.Item(122, 8).Value = 11911625
$ws.Cells.Item(122, 9).Value = 19240172
$ws.Cells.Item(122, 10).Value = 2736.25
$ws.Cells.Item(122, 11).Value = 57720516
$ws.Cells.Item(122, 12).Value = 8208.75
$ws.Cells.Item(122, 13).Value = -57718066
$ws.Cells.Item(122, 14).Value = -13108.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 568.6667
$ws.Cells.Item(81, 9).Value = 568.6667
$ws.Cells.Item(81, 11).Value = 1137.3334
$ws.Cells.Item(81, 13).Value = -76.33339999999998
$ws.Cells.Item(84, 8).Value = 568.6667
$ws.Cells.Item(84, 9).Value = 568.6667
$ws.Cells.Item(84, 11).Value = 5686.666999999999
$ws.Cells.Item(84, 13).Value = -382.6669999999995
$ws.Cells.Item(96, 8).Value = 3365.4167
$ws.Cells.Item(96, 9).Value = 3544.4443
$ws.Cells.Item(96, 11).Value = 3544.4443
$ws.Cells.Item(96, 13).Value = -2171.4443
$ws.Cells.Item(100, 8).Value = 957.0833
$ws.Cells.Item(100, 9).Value = 1130.875
$ws.Cells.Item(100, 11).Value = 2261.75
$ws.Cells.Item(100, 13).Value = -1720.75
$ws.Cells.Item(122, 8).Value = 125001650
$ws.Cells.Item(122, 9).Value = 125001650
$ws.Cells.Item(122, 11).Value = 375004950
$ws.Cells.Item(122, 13).Value = -375002500
$ws.Cells.Item(132, 8).Value = 2836.8235
$ws.Cells.Item(132, 9).Value = 2885.6897
$ws.Cells.Item(132, 11).Value = 8657.069100000001
$ws.Cells.Item(132, 13).Value = -6127.069100000001
$ws.Cells.Item(136, 8).Value = 446.11627
$ws.Cells.Item(136, 9).Value = 264.62964
$ws.Cells.Item(136, 11).Value = 793.88892
$ws.Cells.Item(136, 13).Value = 1756.11108

